$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-17 (Player, Position, Team)
$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
